$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.426.18"
$ws.Range("E2").Value = "  -0.76%  "
$ws.Range("D3").Value = "3.855.48"
$ws.Range("E3").Value = "  -1.93%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "519.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.607"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.31%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.709"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.46%  "
$ws.Range("E10").Value = "  -5.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000317"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -7.75%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "41.39"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.77%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.31"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.38%  "
$ws.Range("D14").Value = "4.475.66"
$ws.Range("E14").Value = "  -1.82%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.40"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.35%  "
$ws.Range("D16").Value = "3.864.08"
$ws.Range("E16").Value = "  -1.62%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.12"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.73%  "
$ws.Range("E18").Value = "  -1.84%  "
$ws.Range("E19").Value = "  +2.40%  "
$ws.Range("D20").Value = "68.509.64"
$ws.Range("E20").Value = "  -0.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "414.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.94%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.97"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "86.55"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.81%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.46"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -8.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.49"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "35.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.89%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "13.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.67%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "677.32"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.59%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.124"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.83%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.73"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +11.12%  "
$ws.Range("B33").Value = "Toncoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.77"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "65.23"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.449"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "39.55"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.06%  "
$ws.Range("B37").Value = "ThetaToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.52"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +14.66%  "
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "0.0₃0828"
$ws.Range("E38").Value = "  -6.83%  "
$ws.Range("E39").Value = "  -1.12%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0473"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.50%  "
$ws.Range("E43").Value = "  +2.93%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.77"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.43%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.39"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.92%  "
$ws.Range("E46").Value = "  -3.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.98"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.89%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.000268"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +12.64%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "143.61"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.38%  "
$ws.Range("E50").Value = "  -3.77%  "
$ws.Range("D51").Value = "0.0₆0336"
$ws.Range("E51").Value = "  -7.79%  "
